$d = $word.ActiveDocument

# --- 1. Heading "Q1 **" -> "Q1" -----------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Q1 **", $true, $false, $false, $false, $false, $true, 1, $false, "Q1", 2)

# --- 2. Q5 answer paragraph: rewrite the ts/pval sentence ----------------
# Replaces the "ts"/"pval" VerbatimChar-styled placeholders with the
# actual computed values (plain text).
$rng = $d.Content
$old2a = "The test statistics is ts. The p-value is pval."
$new2a = "The test statistics is 4.24. The p-value is 4.9594777^{-4}."
$null = $rng.Find.Execute($old2a, $true, $false, $false, $false, $false, $true, 1, $false, $new2a, 2)

$rng = $d.Content
$old2b = "Based on the p-value, XXXXXXXXXXXXX"
$new2b = "The p-value is less than 0.05. Therefore, the null hypothesis should be rejected. We find evidence that there is a difference in means."
$null = $rng.Find.Execute($old2b, $true, $false, $false, $false, $false, $true, 1, $false, $new2b, 2)

# --- 3. Q6 answer paragraph: simplify wording -----------------------------
$rng = $d.Content
$old3 = "The p-value is the probability of observing the test statistic (4.24) as or more supportive of the alternative hypothesis (u1 - u2 does not equal 0) than the actual observed value, given the null hypothesis is true. In other words, the probability;ity that if H0: u1 - u2 = 0 was true, a new sample of data would give a test statistic at least as large as 4.24 in absolute magnitude."
$new3 = "The p-value is the probability that a new sample would produce would produce a test statistic at least as large as 4.24 in magnitude if the null hypothesis were true."
$null = $rng.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# --- 4. "B. T-test output above..." paragraph -----------------------------
$rng = $d.Content
$old4 = "B. T-test output above. The test statistic is ts and the p-value is pval."
$new4 = "B. T-test output above. The test statistic is -0.8609794 and the p-value is 0.4115991."
$null = $rng.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# --- 5. "The p-value (0.41) is greater than alpha..." paragraph ----------
$rng = $d.Content
$old5 = "The p-value (0.41) is greater than alpha, which indicates failing to reject the null hypothesis (difference of means = 0). Therefore, we did not find evidence that there is a difference between means. Additionally, the confidence interval of -9.37 to 4.21 includes 0, indicating that we can conclude that there is not a difference between means. ***"
$new5 = "The p-value (0.41) is greater than alpha, which indicates failing to reject the null hypothesis. Therefore, we did not find evidence that there is a difference between means. Additionally, the confidence interval of -9.37 to 4.21 includes 0, indicating that we can conclude that there is not a difference between means."
$null = $rng.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

# --- 6. Remove the "## do I need to multiple by 2..." comment run --------
$rng = $d.Content
$found = $rng.Find.Execute("## do I need to multiple by 2, how to use to reject or ftr the Ho")
if ($found) {
    $rng.Delete()
}
